$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("L10").ClearContents()
$ws.Range("N10").Value = 0
$ws.Range("H12").Value = 4697.6523
$ws.Range("I12").Value = 5342.85
$ws.Range("J12").Value = 396.33334
$ws.Range("K12").Value = 5342.85
$ws.Range("L12").Value = 396.33334
$ws.Range("M12").Value = -5172.85
$ws.Range("N12").Value = -736.33334
$ws.Range("H17").Value = 126659.25
$ws.Range("J17").Value = 126659.25
$ws.Range("L17").Value = 379977.75
$ws.Range("N17").Value = -380313.75
$ws.Range("H100").Value = 4137.25
$ws.Range("I100").Value = 3000
$ws.Range("K100").Value = 3000
$ws.Range("M100").Value = -2459
$ws.Range("H132").Value = 4940.143
$ws.Range("I132").Value = 4940.143
$ws.Range("K132").Value = 14820.429
$ws.Range("M132").Value = -12290.429
$ws.Range("H138").Value = 3368.2144
$ws.Range("J138").Value = 3699
$ws.Range("L138").Value = 11097
$ws.Range("N138").Value = -21377
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 487.5
$ws.Range("J5").Value = 500
$ws.Range("L5").Value = 500
$ws.Range("N5").Value = -724
$ws.Range("H32").Value = 3685.9312
$ws.Range("I32").Value = 3685.9312
$ws.Range("J32").Value = 0
$ws.Range("K32").Value = 3685.9312
$ws.Range("L32").Value = 0
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -3398.9312
$ws.Range("H36").Value = 10949.667
$ws.Range("I36").Value = 10949.667
$ws.Range("K36").Value = 10949.667
$ws.Range("M36").Value = -10603.667
$ws.Range("H97").Value = 916.625
$ws.Range("I97").Value = 903.6667
$ws.Range("K97").Value = 903.6667
$ws.Range("M97").Value = -407.6667
$ws.Range("H102").Value = 25002802
$ws.Range("I102").Value = 33335404
$ws.Range("K102").Value = 33335404
$ws.Range("M102").Value = -33333782
$ws.Range("H132").Value = 1889934.1
$ws.Range("I132").Value = 2003055.2
$ws.Range("K132").Value = 6009165.6
$ws.Range("M132").Value = -6006635.6
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 487.5
$ws.Range("J4").Value = 500
$ws.Range("L4").Value = 500
$ws.Range("N4").Value = -730
$ws.Range("H35").Value = 9999
$ws.Range("I35").Value = 9999
$ws.Range("K35").Value = 9999
$ws.Range("M35").Value = -9689
$ws.Range("H86").Value = 1994.2941
$ws.Range("I86").Value = 2073.4666
$ws.Range("K86").Value = 2073.4666
$ws.Range("M86").Value = -950.4666000000002
$ws.Range("H89").Value = 1994.2941
$ws.Range("I89").Value = 2073.4666
$ws.Range("K89").Value = 10367.333
$ws.Range("M89").Value = -4751.333000000001
$ws.Range("H99").Value = 1498
$ws.Range("I99").Value = 1498
$ws.Range("J99").Value = 1498
$ws.Range("K99").Value = 1498
$ws.Range("L99").Value = 1498
$ws.Range("M99").Value = 0
$ws.Range("N99").Value = -4494
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 242.5
$ws.Range("J7").Value = 293.44446
$ws.Range("L7").Value = 293.44446
$ws.Range("N7").Value = -519.4444599999999
$ws.Range("H16").Value = 1211657.1
$ws.Range("J16").Value = 6637.25
$ws.Range("L16").Value = 6637.25
$ws.Range("N16").Value = -7211.25
$ws.Range("H22").Value = 17082.5
$ws.Range("I22").Value = 25399
$ws.Range("K22").Value = 25399
$ws.Range("M22").Value = -25049
$ws.Range("H94").Value = 926.5
$ws.Range("I94").Value = 902.6667
$ws.Range("K94").Value = 902.6667
$ws.Range("M94").Value = -451.6667
$ws.Range("H105").Value = 2502894.5
$ws.Range("I105").Value = 3334859.2
$ws.Range("K105").Value = 3334859.2
$ws.Range("M105").Value = -3333112.2
$ws.Range("H107").Value = 1081827.5
$ws.Range("I107").Value = 1812084
$ws.Range("J107").Value = 205519.8
$ws.Range("K107").Value = 1812084
$ws.Range("L107").Value = 205519.8
$ws.Range("M107").Value = -1810164
$ws.Range("N107").Value = -209359.8
$ws.Range("H113").Value = 1211657.1
$ws.Range("J113").Value = 6637.25
$ws.Range("L113").Value = 6637.25
$ws.Range("N113").Value = -10977.25
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 200
$ws.Range("J17").Value = 300
$ws.Range("L17").Value = 900
$ws.Range("N17").Value = -1238
$ws.Range("H76").Value = 14374.875
$ws.Range("I76").Value = 8000
$ws.Range("J76").Value = 18199.8
$ws.Range("K76").Value = 24000
$ws.Range("L76").Value = 54599.39999999999
$ws.Range("M76").Value = -23617
$ws.Range("N76").Value = -55365.39999999999
$ws.Range("H79").Value = 14374.875
$ws.Range("I79").Value = 8000
$ws.Range("J79").Value = 18199.8
$ws.Range("K79").Value = 24000
$ws.Range("L79").Value = 54599.39999999999
$ws.Range("M79").Value = -22674
$ws.Range("N79").Value = -57251.39999999999
$ws.Range("H87").Value = 7302.3335
$ws.Range("I87").Value = 7302.3335
$ws.Range("K87").Value = 21907.0005
$ws.Range("M87").Value = -20659.0005
$ws.Range("H88").Value = 12135.083
$ws.Range("J88").Value = 12135.083
$ws.Range("L88").Value = 36405.249
$ws.Range("N88").Value = -37261.249
$ws.Range("H90").Value = 7302.3335
$ws.Range("I90").Value = 7302.3335
$ws.Range("K90").Value = 65721.0015
$ws.Range("M90").Value = -59481.0015
$ws.Range("H91").Value = 12135.083
$ws.Range("J91").Value = 12135.083
$ws.Range("L91").Value = 36405.249
$ws.Range("N91").Value = -39369.249
$ws.Range("H104").Value = 375
$ws.Range("I104").Value = 375
$ws.Range("K104").Value = 1125
$ws.Range("M104").Value = 1496
$ws.Range("H108").Value = 1901.6
$ws.Range("I108").Value = 1779.6666
$ws.Range("J108").Value = 2999
$ws.Range("K108").Value = 5338.9998
$ws.Range("L108").Value = 8997
$ws.Range("M108").Value = -2458.9998
$ws.Range("N108").Value = -14757
$ws.Range("H111").Value = 5512.75
$ws.Range("I111").Value = 3684
$ws.Range("K111").Value = 11052
$ws.Range("M111").Value = -7985
$ws.Range("H118").Value = 1849.75
$ws.Range("I118").Value = 1849.75
$ws.Range("K118").Value = 5549.25
$ws.Range("M118").Value = -4306.25
$ws.Range("H122").Value = 293.55
$ws.Range("I122").Value = 272.17648
$ws.Range("J122").Value = 414.66666
$ws.Range("K122").Value = 2449.58832
$ws.Range("L122").Value = 3731.99994
$ws.Range("M122").Value = 0.4116799999997056
$ws.Range("N122").Value = -8631.99994
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H41").Value = 18499.666
$ws.Range("I41").Value = 10000
$ws.Range("J41").Value = 22749.5
$ws.Range("K41").Value = 10000
$ws.Range("L41").Value = 22749.5
$ws.Range("M41").Value = -9645
$ws.Range("N41").Value = -23459.5
$ws.Range("H99").Value = 16006.833
$ws.Range("I99").Value = 5208.4
$ws.Range("J99").Value = 69999
$ws.Range("K99").Value = 5208.4
$ws.Range("L99").Value = 69999
$ws.Range("M99").Value = -2962.4
$ws.Range("N99").Value = -74491
$ws.Range("H102").Value = 1599.5
$ws.Range("I102").Value = 1599.5
$ws.Range("K102").Value = 1599.5
$ws.Range("M102").Value = 22.5
$ws.Range("H132").Value = 7816095.5
$ws.Range("I132").Value = 7816095.5
$ws.Range("K132").Value = 23448286.5
$ws.Range("M132").Value = -23445756.5
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H42").Value = 25000
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H49").Value = 25000
$ws.Range("I49").Value = 0
$ws.Range("K49").Value = 0
$ws.Range("M49").ClearContents()
$ws.Range("H122").Value = 6278.5
$ws.Range("I122").Value = 6442.2354
$ws.Range("K122").Value = 19326.7062
$ws.Range("M122").Value = -16876.7062
$ws.Range("H132").Value = 41688652
$ws.Range("J132").Value = 6633.3335
$ws.Range("L132").Value = 19900.0005
$ws.Range("N132").Value = -24960.0005
$ws.Range("H133").Value = 56666.332
$ws.Range("J133").Value = 56666.332
$ws.Range("L133").Value = 56666.332
$ws.Range("N133").Value = -61726.332
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 2988.125
$ws.Range("J96").Value = 2919.6
$ws.Range("L96").Value = 2919.6
$ws.Range("N96").Value = -5665.6
$ws.Range("H100").Value = 1419.8334
$ws.Range("I100").Value = 655
$ws.Range("J100").Value = 2949.5
$ws.Range("K100").Value = 1310
$ws.Range("L100").Value = 5899
$ws.Range("M100").Value = -769
$ws.Range("N100").Value = -6981
$ws.Range("J132").Value = 9450.214
$ws.Range("L132").Value = 28350.642
$ws.Range("N132").Value = -33410.642
